# SVM model Results Compilation
# Append 12 new SVM model result rows (rows 64-75) to Sheet1, matching the
# existing layout: A=S.No., B=Model, C=AUC, D=precision, E=recall,
# F=f1-score, G=training_loss, H=training_accuracy, I=test_loss,
# J=test_accuracy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(61, "Std_SVM_DummyEncoded_MinMaxScaling",       0.88, 0.54, 0.13, 0.21, 3.6774, 0.8935, 3.6638, 0.8939),
    @(62, "Std_SVM_DummyEncoded_MinMaxScaling_SMOTE", 0.85, 0.45, 0.17, 0.25, 2.1255, 0.9384, 3.8928, 0.8873),
    @(63, "Std_SVM_DummyEncoded_SMOTE",               0.84, 0.29, 0.71, 0.42, 8.5347, 0.7529, 7.4804, 0.7834),
    @(64, "Std_SVM_DummyEncoded",                     0.62, 0,    0,    0,    4.0212, 0.8836, 3.7783, 0.8906),
    @(65, "Std_SVM_DummyEncoded_Binning",              0.83, 0.54, 0.34, 0.42, 3.1903, 0.9076, 3.5493, 0.8972),
    @(66, "Std_SVM_DummyEncoded_Binning_SMOTE",        0.69, 0,    0,    0,    3.4885, 0.899,  3.7401, 0.8917),
    @(67, "Std_SVM_LabeEncoded_MinMaxScaling",         0.77, 0,    0,    0,    4.0403, 0.883,  3.7401, 0.8917),
    @(68, "Std_SVM_LabelEncoded_MinMaxScaling_SMOTE",  0.76, 0.23, 0.64, 0.34, 7.1285, 0.7936, 9.5413, 0.7237),
    @(69, "Std_SVM_LabelEncoded_SMOTE",                0.84, 0.3,  0.7,  0.42, 8.5131, 0.7535, 7.3659, 0.7867),
    @(70, "Std_SVM_LabelEncoded",                      0.69, 0,    0,    0,    4.0308, 0.8833, 3.7783, 0.8906),
    @(71, "Std_SVM_LabelEncoded_Binning",              0.7,  0.67, 0.04, 0.08, 3.9448, 0.8858, 3.6638, 0.8939),
    @(72, "Std_SVM_LabelEncoded_Binning_SMOTE",        0.76, 0.4,  0.4,  0.4,  3.6994, 0.8929, 4.4653, 0.8707)
)

$startRow = 64
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
}

# Scroll the viewport so row 50 is the top visible row, then select K76 as
# the active cell (matches the saved view state of the edited workbook).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K76").Select()
